$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 498.22223
$ws.Range("I92").Value = 397.7143
$ws.Range("J92").Value = 850
$ws.Range("K92").Value = 397.7143
$ws.Range("L92").Value = 850
$ws.Range("M92").Value = 850.2857
$ws.Range("N92").Value = -3346

$ws.Range("H98").Value = 34094.87
$ws.Range("I98").Value = 64079.25
$ws.Range("J98").Value = 2111.5334
$ws.Range("K98").Value = 64079.25
$ws.Range("L98").Value = 2111.5334
$ws.Range("M98").Value = -62581.25
$ws.Range("N98").Value = -5107.5334

$ws.Range("H122").Value = 34094.87
$ws.Range("I122").Value = 64079.25
$ws.Range("J122").Value = 2111.5334
$ws.Range("K122").Value = 192237.75
$ws.Range("L122").Value = 6334.600199999999
$ws.Range("M122").Value = -189787.75
$ws.Range("N122").Value = -11234.6002

$ws.Range("H132").Value = 1554305.2
$ws.Range("I132").Value = 2232994.8
$ws.Range("K132").Value = 6698984.399999999
$ws.Range("M132").Value = -6696454.399999999

$ws.Range("H137").Value = 12196409
$ws.Range("I137").Value = 18519344
$ws.Range("K137").Value = 55558032
$ws.Range("M137").Value = -55555482

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 9959.154
$ws.Range("I31").Value = 1963.2222
$ws.Range("J31").Value = 27950
$ws.Range("K31").Value = 1963.2222
$ws.Range("L31").Value = 27950
$ws.Range("M31").Value = -1669.2222
$ws.Range("N31").Value = -28538

$ws.Range("H32").Value = 11042.17
$ws.Range("I32").Value = 6438.57
$ws.Range("K32").Value = 6438.57
$ws.Range("M32").Value = -6151.57

$ws.Range("H70").Value = 82250
$ws.Range("J70").Value = 82250
$ws.Range("L70").Value = 82250
$ws.Range("N70").Value = -82790

$ws.Range("H73").Value = 82250
$ws.Range("J73").Value = 82250
$ws.Range("L73").Value = 82250
$ws.Range("N73").Value = -84122

$ws.Range("H112").Value = 26646.75
$ws.Range("J112").Value = 26646.75
$ws.Range("L112").Value = 26646.75
$ws.Range("N112").Value = -29600.75

$ws.Range("H123").Value = 20766.666
$ws.Range("J123").Value = 20766.666
$ws.Range("L123").Value = 20766.666
$ws.Range("N123").Value = -30566.666

$ws.Range("H132").Value = 1630.6863
$ws.Range("I132").Value = 990.79486
$ws.Range("J132").Value = 3710.3333
$ws.Range("K132").Value = 2972.38458
$ws.Range("L132").Value = 11130.9999
$ws.Range("M132").Value = -442.3845799999999
$ws.Range("N132").Value = -16190.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H27").Value = 73333.336
$ws.Range("J27").Value = 73333.336
$ws.Range("L27").Value = 73333.336
$ws.Range("N27").Value = -73717.336

$ws.Range("H93").Value = 51900
$ws.Range("J93").Value = 51900
$ws.Range("L93").Value = 51900
$ws.Range("N93").Value = -55644

$ws.Range("H97").Value = 35199.8
$ws.Range("I97").Value = 17499.666
$ws.Range("J97").Value = 61750
$ws.Range("K97").Value = 17499.666
$ws.Range("L97").Value = 61750
$ws.Range("M97").Value = -16508.666
$ws.Range("N97").Value = -63732

$ws.Range("H101").Value = 38718
$ws.Range("J101").Value = 38718
$ws.Range("L101").Value = 38718
$ws.Range("N101").Value = -45208

$ws.Range("H105").Value = 1264269.5
$ws.Range("I105").Value = 1624910.9
$ws.Range("J105").Value = 2025
$ws.Range("K105").Value = 1624910.9
$ws.Range("L105").Value = 2025
$ws.Range("M105").Value = -1623163.9
$ws.Range("N105").Value = -5519

$ws.Range("H107").Value = 815.4231
$ws.Range("I107").Value = 838
$ws.Range("J107").Value = 720.6
$ws.Range("K107").Value = 838
$ws.Range("L107").Value = 720.6
$ws.Range("M107").Value = 1082
$ws.Range("N107").Value = -4560.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 66670380
$ws.Range("I105").Value = 166672130
$ws.Range("J105").Value = 2552.6667
$ws.Range("K105").Value = 166672130
$ws.Range("L105").Value = 2552.6667
$ws.Range("M105").Value = -166670383
$ws.Range("N105").Value = -6046.6667

$ws.Range("H106").Value = 28966.666
$ws.Range("J106").Value = 28966.666
$ws.Range("L106").Value = 28966.666
$ws.Range("N106").Value = -31490.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 79.888885
$ws.Range("I12").Value = 19.5
$ws.Range("J12").Value = 110.083336
$ws.Range("K12").Value = 58.5
$ws.Range("L12").Value = 330.250008
$ws.Range("M12").Value = 114.5
$ws.Range("N12").Value = -676.250008

$ws.Range("H95").Value = 9800
$ws.Range("J95").Value = 9800
$ws.Range("L95").Value = 29400
$ws.Range("N95").Value = -33518

$ws.Range("H104").Value = 903.5714
$ws.Range("J104").Value = 1173
$ws.Range("L104").Value = 3519
$ws.Range("N104").Value = -8761

$ws.Range("H131").Value = 11409829
$ws.Range("I131").Value = 35714500
$ws.Range("J131").Value = 67649.8
$ws.Range("K131").Value = 107143500
$ws.Range("L131").Value = 202949.4
$ws.Range("M131").Value = -107138460
$ws.Range("N131").Value = -213029.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 52294532
$ws.Range("I11").Value = 5600000
$ws.Range("J11").Value = 84621510
$ws.Range("K11").Value = 5600000
$ws.Range("L11").Value = 84621510
$ws.Range("M11").Value = -5599861
$ws.Range("N11").Value = -84621788

$ws.Range("H135").Value = 39422
$ws.Range("J135").Value = 39422
$ws.Range("L135").Value = 39422
$ws.Range("N135").Value = -49562

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1000
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1000
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1000
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -1590

$ws.Range("H27").Value = 1000
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1000
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1000
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1214

$ws.Range("H46").Value = 920.8
$ws.Range("I46").Value = 886.8570999999999
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 886.8570999999999
$ws.Range("L46").Value = 1000
$ws.Range("M46").Value = -698.8570999999999
$ws.Range("N46").Value = -1376

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 20562.5
$ws.Range("J28").Value = 20562.5
$ws.Range("L28").Value = 20562.5
$ws.Range("N28").Value = -21258.5

$ws.Range("H70").Value = 27113.125
$ws.Range("J70").Value = 28129.285
$ws.Range("L70").Value = 28129.285
$ws.Range("N70").Value = -28759.285

$ws.Range("H73").Value = 27113.125
$ws.Range("J73").Value = 28129.285
$ws.Range("L73").Value = 28129.285
$ws.Range("N73").Value = -30313.285

$ws.Range("H119").Value = 9349
$ws.Range("J119").Value = 9349
$ws.Range("L119").Value = 9349
$ws.Range("N119").Value = -19025

Write-Host "Applied all changes"